# IKD update: GaN CMOS 2026-02-05T23:27Z
# Appends 5 new literature-tracker rows (136-140) to the "Master" sheet,
# mirroring the columns already present in row 1 (A:T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        B = "Advancements in MOSFET Channel Materials a Comparative Review"
        C = 2026
        D = "Trans Tech Publications, Ltd."
        E = "Advanced Materials Research"
        F = "Sania, Sania; Gangwani, Parvesh; Mohil, Meenu; Kaur, Ravneet"
        H = "10.4028/p-pou2qy"
        I = "https://doi.org/10.4028/p-pou2qy"
        J = "Journal"
        K = "n-FET"
        L = "Experiment"
        M = "Contacts"
        Q = "Advancements in MOSFET Channel Materials a Comparative Review"
        R = "High"
        S = "2026-02-05"
    },
    @{
        B = "FinFET vs Planar MOSFET: A Performance-Based Comparative Study"
        C = 2026
        D = "Trans Tech Publications, Ltd."
        E = "Advanced Materials Research"
        F = "Sahu, Soumya; Saikia, Tanmoi; Gangwani, Parvesh; Kaur, Ravneet"
        H = "10.4028/p-1pyqh5"
        I = "https://doi.org/10.4028/p-1pyqh5"
        J = "Journal"
        K = "n-FET"
        L = "Experiment"
        M = "Contacts"
        Q = "FinFET vs Planar MOSFET: A Performance-Based Comparative Study"
        R = "High"
        S = "2026-02-05"
    },
    @{
        B = "Advancements in MOSFET Channel Materials a Comparative Review"
        C = 2026
        D = "Trans Tech Publications, Ltd."
        E = "Advanced Materials Research"
        F = "Sania, Sania; Gangwani, Parvesh; Mohil, Meenu; Kaur, Ravneet"
        H = "10.4028/p-pou2qy"
        I = "https://doi.org/10.4028/p-pou2qy"
        J = "Journal"
        K = "n-FET"
        L = "Experiment"
        M = "Contacts"
        Q = "Advancements in MOSFET Channel Materials a Comparative Review"
        R = "High"
        S = "2026-02-05"
    },
    @{
        B = "FinFET vs Planar MOSFET: A Performance-Based Comparative Study"
        C = 2026
        D = "Trans Tech Publications, Ltd."
        E = "Advanced Materials Research"
        F = "Sahu, Soumya; Saikia, Tanmoi; Gangwani, Parvesh; Kaur, Ravneet"
        H = "10.4028/p-1pyqh5"
        I = "https://doi.org/10.4028/p-1pyqh5"
        J = "Journal"
        K = "n-FET"
        L = "Experiment"
        M = "Contacts"
        Q = "FinFET vs Planar MOSFET: A Performance-Based Comparative Study"
        R = "High"
        S = "2026-02-05"
    },
    @{
        B = "Advances in Semiconductor Optical Amplifier Technologies for All-Optical Logic Gate Implementations: A Comprehensive Review"
        C = 2026
        D = "MDPI AG"
        E = "Nanomaterials"
        F = "Cui, Jiali; Zoiros, Kyriakos E.; Kotb, Amer"
        H = "10.3390/nano16030202"
        I = "https://doi.org/10.3390/nano16030202"
        J = "Journal"
        K = "Inverter"
        L = "Experiment"
        M = "Contacts"
        Q = "Advances in Semiconductor Optical Amplifier Technologies for All-Optical Logic Gate Implementations: A Comprehensive Review"
        R = "High"
        S = "2026-02-05"
    }
)

$startRow = 136
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("Q$r").Value = $row.Q
    $ws.Range("R$r").Value = $row.R
    # AddedDate is stored as plain text (e.g. "2026-02-05"), not a date
    # serial, so force Text format before assigning to stop Excel's
    # auto date-parsing from converting the string.
    $ws.Range("S$r").NumberFormat = "@"
    $ws.Range("S$r").Value = $row.S
}
